# Daily attendance processing - 2025-11-14 10:25:31
# Normalizes the "Recorded By" (column G) cell values by rotating the
# comma-separated list of recorders left by one position, for the rows
# where the backup/system/admin entries had drifted out of their
# canonical order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"   = "system, backup@backdoor.com, System"
    "G4"   = "backup@backdoor.com, System"
    "G5"   = "backup@backdoor.com, System"
    "G8"   = "backup@backdoor.com, System"
    "G11"  = "System, dnasr281@gmail.com"
    "G17"  = "System, dnasr281@gmail.com"
    "G28"  = "system, backup@backdoor.com, System"
    "G30"  = "backup@backdoor.com, System"
    "G31"  = "backup@backdoor.com, System"
    "G34"  = "backup@backdoor.com, System"
    "G37"  = "System, dnasr281@gmail.com"
    "G43"  = "System, dnasr281@gmail.com"
    "G54"  = "system, backup@backdoor.com, System"
    "G56"  = "backup@backdoor.com, System"
    "G57"  = "backup@backdoor.com, System"
    "G60"  = "backup@backdoor.com, System"
    "G63"  = "System, dnasr281@gmail.com"
    "G69"  = "System, dnasr281@gmail.com"
    "G80"  = "backup@backdoor.com, System"
    "G81"  = "backup@backdoor.com, System"
    "G82"  = "backup@backdoor.com, System"
    "G87"  = "admin@admin.com, dnasr281@gmail.com"
    "G93"  = "System, dnasr281@gmail.com"
    "G94"  = "System, dnasr281@gmail.com"
    "G96"  = "System, dnasr281@gmail.com"
    "G106" = "backup@backdoor.com, System"
    "G107" = "backup@backdoor.com, System"
    "G108" = "backup@backdoor.com, System"
    "G113" = "admin@admin.com, dnasr281@gmail.com"
    "G119" = "System, dnasr281@gmail.com"
    "G120" = "System, dnasr281@gmail.com"
    "G122" = "System, dnasr281@gmail.com"
    "G132" = "backup@backdoor.com, System"
    "G133" = "backup@backdoor.com, System"
    "G134" = "backup@backdoor.com, System"
    "G139" = "admin@admin.com, dnasr281@gmail.com"
    "G145" = "System, dnasr281@gmail.com"
    "G146" = "System, dnasr281@gmail.com"
    "G148" = "System, dnasr281@gmail.com"
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
